$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Input" sheet (tabSelected)

# Row 3 (A3:B3 = "active"/"checked") already carries the label/value
# formatting (green-filled label cell + wrap-text value cell) that the
# new row should reuse, so copy its formatting down onto row 6 first.
$ws.Range("A3:B3").Copy($ws.Range("A6:B6"))

# Fill in the new submit-client row.
$ws.Range("A6").Value = "submitclient"
$ws.Range("B6").Value = "submit"

$ws.Range("A6").Select()
